$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the year table by one column (K) for the new 2023 data ---
# Mirror the existing last column's (J) number/border formatting onto the
# new column so the new cells inherit the same fonts / number format as
# the rest of the table.
$ws.Range("J3:J6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- New 2023 values ---
$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 7562
$ws.Range("K5").Value = 1514
$ws.Range("K6").Value = 6048

# Column K is now the right-most column of the table, so it needs its own
# right-hand border closing off the box (matching the border already used
# on the left/around the rest of the table).
$ws.Range("K3:K6").Borders.Item(10).LineStyle = 1   # xlEdgeRight / xlContinuous

# Widen the following empty columns (L:O) to match the table's standard
# column width, same as the rest of the data columns (B:J), anticipating
# future years being added.
$tableColWidth = $ws.Columns("J").ColumnWidth
$ws.Columns("K:O").ColumnWidth = $tableColWidth

$excel.CutCopyMode = 0
